$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cash_Spent")

# ---------------------------------------------------------------------------
# 1) Two new expense rows were typed in at the top of the Cash Spent table.
#    That pushes every existing data row (old row 6 .. old row 17) down by two
#    rows. Shift bottom-up (columns B:I, which also carries the stray I7
#    helper cell along for the ride) so we never clobber a source row before
#    it has been read.
# ---------------------------------------------------------------------------
for ($r = 17; $r -ge 6; $r--) {
    $target = $r + 2
    $ws.Range("B$target`:I$target").Formula = $ws.Range("B$r`:I$r").Formula
}

# ---------------------------------------------------------------------------
# 2) The old row 5 ("Cine" / 80000) slides down to row 7, picking up the
#    banded-row style (30) the table applies to freshly-touched rows, and
#    gains its own B7 date cell.
# ---------------------------------------------------------------------------
$ws.Range("B7").Value = " 31/10/2022"
$ws.Range("C7").Value = $ws.Range("C5").Value2
$ws.Range("D7").Value = $ws.Range("D5").Value2
$ws.Range("E7").Value = $ws.Range("E5").Value2
$ws.Range("F7").Value = $ws.Range("F5").Value2
$ws.Range("B7:F7").Style = "Normal"

# ---------------------------------------------------------------------------
# 3) Row 6 is a brand-new entry: "Cine" 50000.
# ---------------------------------------------------------------------------
$ws.Range("B6").Value = " 31/10/2022"
$ws.Range("C6").Value = "21:30:56"
$ws.Range("D6").Value = """Cine"
$ws.Range("E6").Value = "50000"
$ws.Range("F6").Value = "Gasto"
$ws.Range("B6:F6").Style = "Normal"

# ---------------------------------------------------------------------------
# 4) Row 5 is also a brand-new entry: "Arriendo" 800000 (the date in B5 is
#    left as-is).
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "21:31:04"
$ws.Range("D5").Value = """Arriendo"
$ws.Range("E5").Value = "800000"
$ws.Range("F5").Value = "Gasto"

# ---------------------------------------------------------------------------
# 5) The two slicer "helper" notes (J9 / J13) are lost in the process — clear
#    them out while keeping the merged ranges (J8:K11, J12:K19) intact.
# ---------------------------------------------------------------------------
$ws.Range("J8:K11").UnMerge()
$ws.Range("J8:K11").Merge()
$ws.Range("J12:K19").UnMerge()
$ws.Range("J12:K19").Merge()

# ---------------------------------------------------------------------------
# 6) Touch up the chart axis / legend fonts on both charts — Excel leaves a
#    stray empty run behind in the txPr the first time you do this, and this
#    edit fills that run in with literal text "None".
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Cash Summary")
foreach ($chartObj in $ws1.ChartObjects()) {
    $chart = $chartObj.Chart
    $chart.Axes(2).TickLabels.Font.Name = $chart.Axes(2).TickLabels.Font.Name
}

$ws4 = $wb.Worksheets.Item("Chart Data")
foreach ($chartObj in $ws4.ChartObjects()) {
    $chart = $chartObj.Chart
    $chart.Axes(1).TickLabels.Font.Name = $chart.Axes(1).TickLabels.Font.Name
    $chart.Axes(2).TickLabels.Font.Name = $chart.Axes(2).TickLabels.Font.Name
    $chart.Legend.Font.Name = $chart.Legend.Font.Name
}
